$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.467.82'
$ws.Range("E2").Value = '  +2.20%  '
$ws.Range("D3").Value = '2.488.80'
$ws.Range("E3").Value = '  +3.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.56'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.81'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("D9").Value = '2.487.77'
$ws.Range("E9").Value = '  +3.46%  '
$ws.Range("E10").Value = '  +1.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.75'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.358'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.68'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +7.43%  '
$ws.Range("D15").Value = '2.926.95'
$ws.Range("E15").Value = '  +3.26%  '
$ws.Range("D16").Value = '63.248.54'
$ws.Range("E16").Value = '  +1.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000143'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.74%  '
$ws.Range("D18").Value = '2.473.79'
$ws.Range("E18").Value = '  +3.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.34'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '342.15'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.42%  '
$ws.Range("E21").Value = '  +1.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.78'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.64%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.68'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.94%  '
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.52'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("E28").Value = '  +4.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.12'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.63%  '
$ws.Range("E30").Value = '  +3.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.78'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +6.08%  '
$ws.Range("D32").Value = '0.0₃0805'
$ws.Range("E32").Value = '  +4.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '176.41'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.76%  '
$ws.Range("E34").Value = '  +8.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '408.50'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +13.46%  '
$ws.Range("E36").Value = '  +1.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.93'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.26%  '
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.37'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.39%  '
$ws.Range("E40").Value = '  +4.43%  '
$ws.Range("E41").Value = '  -0.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.54'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '150.51'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.73'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.78'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.603'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0966'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0521'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0230'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.23'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.35%  '
$ws.Range("D51").Value = '0.0₆0230'
$ws.Range("E51").Value = '  +5.76%  '
